$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1; Col=1; Old="40+6=46"; New="15+29=44"},
    @{Row=1; Col=2; Old="97-51=46"; New="19+68=87"},
    @{Row=1; Col=3; Old="45+50=95"; New="19+19=38"},
    @{Row=1; Col=4; Old="98-95=3"; New="29+0=29"},
    @{Row=1; Col=5; Old="61-52=9"; New="14+40=54"},
    @{Row=2; Col=1; Old="94-47=47"; New="85-5=80"},
    @{Row=2; Col=2; Old="51+31=82"; New="99-80=19"},
    @{Row=2; Col=3; Old="15+38=53"; New="51-34=17"},
    @{Row=2; Col=4; Old="44+5=49"; New="84-32=52"},
    @{Row=2; Col=5; Old="19-13=6"; New="76-30=46"},
    @{Row=3; Col=1; Old="76-73=3"; New="97-75=22"},
    @{Row=3; Col=2; Old="34-17=17"; New="94-74=20"},
    @{Row=3; Col=3; Old="91-1=90"; New="55-16=39"},
    @{Row=3; Col=4; Old="88-34=54"; New="13+8=21"},
    @{Row=3; Col=5; Old="4+62=66"; New="18+56=74"},
    @{Row=4; Col=1; Old="55-39=16"; New="34+36=70"},
    @{Row=4; Col=2; Old="62+29=91"; New="83-40=43"},
    @{Row=4; Col=3; Old="14+10=24"; New="48+43=91"},
    @{Row=4; Col=4; Old="56+34=90"; New="15-0=15"},
    @{Row=4; Col=5; Old="67-45=22"; New="10+2=12"},
    @{Row=5; Col=1; Old="99-5=94"; New="5+72=77"},
    @{Row=5; Col=2; Old="70+7=77"; New="40+27=67"},
    @{Row=5; Col=3; Old="35+1=36"; New="94+4=98"},
    @{Row=5; Col=4; Old="33+13=46"; New="84-41=43"},
    @{Row=5; Col=5; Old="66+33=99"; New="41-18=23"},
    @{Row=6; Col=1; Old="64-17=47"; New="63-59=4"},
    @{Row=6; Col=2; Old="55+15=70"; New="32+23=55"},
    @{Row=6; Col=3; Old="64-31=33"; New="22+56=78"},
    @{Row=6; Col=4; Old="6+45=51"; New="62-32=30"},
    @{Row=6; Col=5; Old="55-3=52"; New="65-20=45"},
    @{Row=7; Col=1; Old="81-7=74"; New="46+46=92"},
    @{Row=7; Col=2; Old="63-31=32"; New="57+25=82"},
    @{Row=7; Col=3; Old="10+13=23"; New="19+62=81"},
    @{Row=7; Col=4; Old="85-54=31"; New="1+67=68"},
    @{Row=7; Col=5; Old="18+42=60"; New="21+17=38"},
    @{Row=8; Col=1; Old="30+7=37"; New="6+61=67"},
    @{Row=8; Col=2; Old="11+43=54"; New="85+5=90"},
    @{Row=8; Col=3; Old="40+7=47"; New="39+38=77"},
    @{Row=8; Col=4; Old="2+80=82"; New="47+15=62"},
    @{Row=8; Col=5; Old="68-60=8"; New="49-45=4"},
    @{Row=9; Col=1; Old="42-3=39"; New="57-15=42"},
    @{Row=9; Col=2; Old="65+21=86"; New="40+3=43"},
    @{Row=9; Col=3; Old="64-15=49"; New="97-7=90"},
    @{Row=9; Col=4; Old="7+53=60"; New="73-69=4"},
    @{Row=9; Col=5; Old="49+33=82"; New="87-10=77"},
    @{Row=10; Col=1; Old="19+2=21"; New="89-45=44"},
    @{Row=10; Col=2; Old="24+50=74"; New="44-27=17"},
    @{Row=10; Col=3; Old="79-21=58"; New="82-10=72"},
    @{Row=10; Col=4; Old="13+54=67"; New="23+49=72"},
    @{Row=10; Col=5; Old="73-52=21"; New="71-45=26"},
    @{Row=11; Col=1; Old="47-29=18"; New="3+84=87"},
    @{Row=11; Col=2; Old="28+59=87"; New="25-18=7"},
    @{Row=11; Col=3; Old="56-25=31"; New="88-69=19"},
    @{Row=11; Col=4; Old="48+12=60"; New="55-36=19"},
    @{Row=11; Col=5; Old="99-18=81"; New="65+30=95"},
    @{Row=12; Col=1; Old="25+74=99"; New="2+84=86"},
    @{Row=12; Col=2; Old="28+12=40"; New="35-16=19"},
    @{Row=12; Col=3; Old="39+39=78"; New="23+18=41"},
    @{Row=12; Col=4; Old="73-59=14"; New="49+34=83"},
    @{Row=12; Col=5; Old="0+90=90"; New="80-48=32"},
    @{Row=13; Col=1; Old="97-41=56"; New="38+50=88"},
    @{Row=13; Col=2; Old="33+45=78"; New="76-20=56"},
    @{Row=13; Col=3; Old="43-10=33"; New="42-33=9"},
    @{Row=13; Col=4; Old="75-75=0"; New="36+17=53"},
    @{Row=13; Col=5; Old="46+4=50"; New="97-42=55"},
    @{Row=14; Col=1; Old="42-3=39"; New="76-1=75"},
    @{Row=14; Col=2; Old="23+45=68"; New="50+31=81"},
    @{Row=14; Col=3; Old="2+44=46"; New="20+19=39"},
    @{Row=14; Col=4; Old="2+79=81"; New="77+7=84"},
    @{Row=14; Col=5; Old="62+10=72"; New="72-17=55"},
    @{Row=15; Col=1; Old="93-9=84"; New="56+18=74"},
    @{Row=15; Col=2; Old="83-34=49"; New="86+1=87"},
    @{Row=15; Col=3; Old="91-41=50"; New="67-32=35"},
    @{Row=15; Col=4; Old="65-54=11"; New="41+49=90"},
    @{Row=15; Col=5; Old="55-32=23"; New="54+32=86"},
    @{Row=16; Col=1; Old="8+38=46"; New="71-67=4"},
    @{Row=16; Col=2; Old="4+64=68"; New="24-7=17"},
    @{Row=16; Col=3; Old="83-55=28"; New="3+32=35"},
    @{Row=16; Col=4; Old="4+57=61"; New="26+7=33"},
    @{Row=16; Col=5; Old="18+6=24"; New="44+4=48"},
    @{Row=17; Col=1; Old="20-2=18"; New="56-51=5"},
    @{Row=17; Col=2; Old="20+64=84"; New="21-20=1"},
    @{Row=17; Col=3; Old="38+37=75"; New="91-50=41"},
    @{Row=17; Col=4; Old="59-15=44"; New="95-24=71"},
    @{Row=17; Col=5; Old="82-30=52"; New="38-13=25"},
    @{Row=18; Col=1; Old="0+12=12"; New="36+31=67"},
    @{Row=18; Col=2; Old="17+0=17"; New="87-2=85"},
    @{Row=18; Col=3; Old="86-20=66"; New="11-5=6"},
    @{Row=18; Col=4; Old="90-24=66"; New="96-94=2"},
    @{Row=18; Col=5; Old="20+66=86"; New="78-47=31"},
    @{Row=19; Col=1; Old="60+21=81"; New="72-54=18"},
    @{Row=19; Col=2; Old="22+54=76"; New="88-0=88"},
    @{Row=19; Col=3; Old="54+17=71"; New="85-1=84"},
    @{Row=19; Col=4; Old="26+45=71"; New="82-28=54"},
    @{Row=19; Col=5; Old="84-39=45"; New="36+34=70"},
    @{Row=20; Col=1; Old="59+35=94"; New="90-60=30"},
    @{Row=20; Col=2; Old="87-82=5"; New="19+62=81"},
    @{Row=20; Col=3; Old="86-19=67"; New="18+26=44"},
    @{Row=20; Col=4; Old="98-83=15"; New="39+52=91"},
    @{Row=20; Col=5; Old="59-35=24"; New="28+71=99"}
)

$mismatchCount = 0
foreach ($r in $replacements) {
    $cell = $t.Cell($r.Row, $r.Col)
    $range = $cell.Range
    $range.MoveEnd(1, -1) | Out-Null
    $current = $range.Text
    if ($current -ne $r.Old) {
        $mismatchCount = $mismatchCount + 1
        Write-Output ("MISMATCH at Row=" + $r.Row + " Col=" + $r.Col + " expected=" + $r.Old + " actual=" + $current)
    }
    $range.Text = $r.New
}
Write-Output ("Done. Mismatches=" + $mismatchCount)
